# Adicionado login e lógica de validade do plano
# Add a new "Validade" column (G) to the quote table, indicating the
# validity month of each insurer's pricing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header -----------------------------------------------------------
$ws.Range("G1").Value = "Validade"

# Match the formatting already used for the "Abrangência" header (F1):
# bold font, thin left/right border, centered/top aligned.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Data ---------------------------------------------------------------
# Hapvida block (rows 2-21) is valid through 2025-12
$ws.Range("G2:G21").Value = "2025-12"

# Samel block (rows 22-41) is valid through 2025-08
$ws.Range("G22:G41").Value = "2025-08"

# Adventist block (rows 42-71) is valid through 2025-12
$ws.Range("G42:G71").Value = "2025-12"

Write-Host "Validade column added"
